$d = $word.ActiveDocument

$script:splitCounter = 0
function Force-RunSplit($doc, $pos) {
    # Inserting then deleting a zero-width bookmark forces the OOXML writer to keep
    # the text on either side of $pos in separate <w:r> runs, mimicking how Word
    # splits runs around an edit point even after the marker itself is gone.
    $script:splitCounter = $script:splitCounter + 1
    $markName = "TMP_SPLIT_" + $script:splitCounter
    $r = $doc.Range($pos, $pos)
    $doc.Bookmarks.Add($markName, $r)
    $doc.Bookmarks.Item($markName).Delete()
}

# =========================================================================
# Change 1: "...performances conducted to form the dataset..."
#        -> "...performances conducted on each letter to form the dataset..."
# with the _GoBack bookmark left sitting right after "on each".
# =========================================================================
$find = $d.Content.Find
$found = $find.Execute("conducted to form the dataset", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$ip = $find.Parent.Duplicate
$ip.Collapse(1)
$ip.MoveStart(1, 10)            # skip past "conducted "
$startPos = $ip.Start
$ip.InsertBefore("on each letter ")

Force-RunSplit $d $startPos                 # split "conducted " | "on each"

$goBackPos = $startPos + 7                   # right after "on each"
$goBackRange = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $goBackRange)    # zero-length _GoBack bookmark

$afterLetterPos = $startPos + 15
Force-RunSplit $d $afterLetterPos            # split " letter " | "to form..."

# =========================================================================
# Change 2: "...With this information in mind, the two-sample t-testing method..."
#        -> "...With this information in mind, the paired t-testing method..."
# =========================================================================
$find2 = $d.Content.Find
$found2 = $find2.Execute("two-sample", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng2 = $find2.Parent.Duplicate
$pos2 = $rng2.Start
$rng2.Delete()
$rng2.InsertBefore("paired")

Force-RunSplit $d $pos2                      # split "...the " | "paired"
$pos2b = $pos2 + 6
Force-RunSplit $d $pos2b                     # split "paired" | " t-testing..."

# =========================================================================
# Change 3: "...individually..." / "...of fifty random observations." runs
# simply get merged with their preceding space run (no visible text change).
# =========================================================================
$find3 = $d.Content.Find
$found3 = $find3.Execute(" individually", $true, $false, $false, $false, $false, $true, 1, $false, " individually", 2)

$find4 = $d.Content.Find
$found4 = $find4.Execute(" of fifty random observations.", $true, $false, $false, $false, $false, $true, 1, $false, " of fifty random observations.", 2)

# =========================================================================
# Change 4: "Using the t-testing function built for two independent samples
# of scores from SciPy's statistical module..."
#        -> "Using the t-testing function built for paired samples from
# SciPy's statistical module..."
# =========================================================================
$find5 = $d.Content.Find
$found5 = $find5.Execute("two independent samples of scores from", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng5 = $find5.Parent.Duplicate
$pos5 = $rng5.Start
$rng5.Delete()
$rng5.InsertBefore("paired samples from")

Force-RunSplit $d $pos5                      # split "built for " | "paired"
$pos5b = $pos5 + 6
Force-RunSplit $d $pos5b                     # split "paired" | " samples from..."
